# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig -- updates metadata
# (Version, Date, Publisher, Jurisdiction) and element Short/Definition text.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Remove the duplicated "Contact" / "No display for ContactDetail" row
# (originally row 11) so everything below shifts up by one row.
$meta.Rows.Item(11).Delete()

# Update Version
$meta.Range("B3").Value = "6.0.0"

# Update Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# Replace the old "Contact" row with a new "Jurisdiction" row
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet "Elements" ---
$elem = $wb.Worksheets.Item("Elements")

# Update Short / Definition text for the root Extension element (row 2)
$elem.Range("K2").Value = "Employee Birth Date"
$elem.Range("L2").Value = "Birthdate of the employee or contract holder. For security, this field may be hidden or the value stored in this field should only contain the year"
